$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old CKD Stage sub-rows (row 26 = "1.0", row 27 = "2.0") are removed.
# Deleting these rows shifts everything below up by 2 (old row 28 "3.0" becomes
# new row 26, old row 29 "Diabetes Type..." becomes new row 27, etc.) and Excel
# automatically adjusts the merged cell ranges (A25:A28 -> A25:A26, A29:A31 -> A27:A29).
$ws.Rows("26:27").Delete()

# Update the CKD Stage, n (%) summary row (row 25) with new values.
# Use a leading apostrophe so Excel stores these numeric-looking values as text,
# matching the other categorical "n (%)" columns in this sheet.
$ws.Range("B25").Value = "'0"
$ws.Range("C25").Value = "1379 (91.7)"
$ws.Range("D25").Value = "6680 (94.4)"

# Update the remaining CKD Stage sub-row (now row 26, previously row 28 "3.0").
$ws.Range("B26").Value = "'3"
